# Daily odds-feed refresh for "Netherlands Eredivisie" sheet
# (commit: "Atualizacao de bases das ligas, do dia: 02-05-2024 as 20:28")
#
# The upstream scraper re-pulled match rows 279-292. Two pairs of rows
# (279/280 and 284/285) came back in swapped order versus the previous
# pull, and a handful of other rows simply got refreshed odds. This
# script writes every affected cell to its new value directly.
#
# Note: Value2 is used throughout because reading/writing via the plain
# .Value property surfaces a proxy object in this host rather than the
# underlying scalar; .Value2 round-trips cleanly for both numbers and text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Rows 279 / 280 (ids 277 / 278): the two fixtures swapped places.
# Column A (row index) stays put; everything B:AB trades between rows.
# ---------------------------------------------------------------------
$ws.Range("B279").Value2 = 7061102
$ws.Range("E279").Value2 = "PEC Zwolle"
$ws.Range("F279").Value2 = "Heracles"
$ws.Range("G279").Value2 = 3
$ws.Range("H279").Value2 = 1
$ws.Range("I279").Value2 = "H"
$ws.Range("J279").Value2 = 2.25
$ws.Range("K279").Value2 = 3.75
$ws.Range("L279").Value2 = 2.75
$ws.Range("M279").Value2 = 2.05
$ws.Range("O279").Value2 = 3.25
$ws.Range("P279").Value2 = -0.5
$ws.Range("Q279").Value2 = 2.025
$ws.Range("R279").Value2 = 1.825
$ws.Range("S279").Value2 = 3
$ws.Range("T279").Value2 = 1.825
$ws.Range("U279").Value2 = 2.025
$ws.Range("V279").Value2 = 1.05
$ws.Range("X279").Value2 = -1
$ws.Range("Y279").Value2 = 1.025
$ws.Range("Z279").Value2 = -1
$ws.Range("AA279").Value2 = 0.825
$ws.Range("AB279").Value2 = -1

$ws.Range("B280").Value2 = 7062784
$ws.Range("E280").Value2 = "NEC"
$ws.Range("F280").Value2 = "AZ"
$ws.Range("G280").Value2 = 0
$ws.Range("H280").Value2 = 3
$ws.Range("I280").Value2 = "A"
$ws.Range("J280").Value2 = 3.2
$ws.Range("K280").Value2 = 4
$ws.Range("L280").Value2 = 1.95
$ws.Range("M280").Value2 = 3.1
$ws.Range("O280").Value2 = 2.1
$ws.Range("P280").Value2 = 0.25
$ws.Range("Q280").Value2 = 2.03
$ws.Range("R280").Value2 = 1.87
$ws.Range("S280").Value2 = 2.75
$ws.Range("T280").Value2 = 1.925
$ws.Range("U280").Value2 = 1.925
$ws.Range("V280").Value2 = -1
$ws.Range("X280").Value2 = 1.1
$ws.Range("Y280").Value2 = -1
$ws.Range("Z280").Value2 = 0.8700000000000001
$ws.Range("AA280").Value2 = 0.4625
$ws.Range("AB280").Value2 = -0.5

# ---------------------------------------------------------------------
# Rows 284 / 285 (ids 282 / 283): fixtures swapped too, including their
# scraper "id" (B col). Those ids look numeric ("7093668"/"7093671")
# but are stored as text, so force text entry (leading apostrophe) and
# then strip the resulting number-format back to Normal/General so no
# stray style gets attached to the cell.
# ---------------------------------------------------------------------
$c = $ws.Range("B284")
$c.Value2 = "'7093671"
$c.Style = "Normal"
$ws.Range("E284").Value2 = "Fortuna Sittard"
$ws.Range("F284").Value2 = "Go Ahead Eagles"
$ws.Range("J284").Value2 = 2.375
$ws.Range("K284").Value2 = 3.5
$ws.Range("L284").Value2 = 2.875
$ws.Range("M284").Value2 = 2.3
$ws.Range("N284").Value2 = 3.5
$ws.Range("O284").Value2 = 3
$ws.Range("P284").Value2 = -0.25
$ws.Range("Q284").Value2 = 2.05
$ws.Range("R284").Value2 = 1.85
$ws.Range("S284").Value2 = 2.5
$ws.Range("T284").Value2 = 1.925
$ws.Range("U284").Value2 = 1.925

$c = $ws.Range("B285")
$c.Value2 = "'7093668"
$c.Style = "Normal"
$ws.Range("E285").Value2 = "Almere City FC"
$ws.Range("F285").Value2 = "Heerenveen"
$ws.Range("J285").Value2 = 2.7
$ws.Range("K285").Value2 = 3.4
$ws.Range("L285").Value2 = 2.55
$ws.Range("M285").Value2 = 2.8
$ws.Range("N285").Value2 = 3.4
$ws.Range("O285").Value2 = 2.5
$ws.Range("P285").Value2 = 0
$ws.Range("Q285").Value2 = 2.07
$ws.Range("R285").Value2 = 1.83
$ws.Range("S285").Value2 = 2.75
$ws.Range("T285").Value2 = 2
$ws.Range("U285").Value2 = 1.85

# ---------------------------------------------------------------------
# Rows 286-292: same fixtures, just refreshed odds (no row reordering).
# ---------------------------------------------------------------------
# Row 286 (id 284, Heracles v RKC)
$ws.Range("M286").Value2 = 2.3
$ws.Range("O286").Value2 = 3
$ws.Range("Q286").Value2 = 2.01
$ws.Range("R286").Value2 = 1.89

# Row 287 (id 285, PSV v Sparta Rotterdam)
$ws.Range("Q287").Value2 = 2.05
$ws.Range("R287").Value2 = 1.85

# Row 288 (id 286, FC Utrecht v Vitesse)
$ws.Range("M288").Value2 = 1.363
$ws.Range("N288").Value2 = 5.25
$ws.Range("O288").Value2 = 8
$ws.Range("P288").Value2 = -1.5
$ws.Range("Q288").Value2 = 2.02
$ws.Range("R288").Value2 = 1.88
$ws.Range("T288").Value2 = 1.8
$ws.Range("U288").Value2 = 2.05

# Row 289 (id 287, FC Volendam v Ajax)
$ws.Range("M289").Value2 = 8
$ws.Range("N289").Value2 = 5.5
$ws.Range("O289").Value2 = 1.363
$ws.Range("Q289").Value2 = 2.01
$ws.Range("R289").Value2 = 1.89

# Row 290 (id 288, AZ v FC Twente)
$ws.Range("M290").Value2 = 2.3
$ws.Range("O290").Value2 = 3.1

# Row 291 (id 289, Feyenoord v PEC Zwolle)
$ws.Range("O291").Value2 = 21
$ws.Range("Q291").Value2 = 1.91
$ws.Range("R291").Value2 = 1.99
$ws.Range("T291").Value2 = 1.825
$ws.Range("U291").Value2 = 2.025

# Row 292 (id 290, Excelsior v NEC)
$ws.Range("N292").Value2 = 3.6
$ws.Range("O292").Value2 = 2.3
$ws.Range("Q292").Value2 = 1.91
$ws.Range("R292").Value2 = 1.99
